$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that immediately follows
#    the H1 title paragraph at the top of the document.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph (re-using the title text) right
#    before the final "Prompt: ..." paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$beforeLastPara = $d.Paragraphs.Item($count - 1)
$insertPoint = $beforeLastPara.Range.End - 1
$insertRange = $d.Range($insertPoint, $insertPoint)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aztec Warrior for Free - A Review of the Slot Game</w:t></w:r></w:p>'
[void]$insertRange.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3. Replace the text of the (now) final paragraph - the old
#    "Prompt: ..." image-generation text - with the meta description
#    text, keeping its italic run formatting untouched.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
[void]$lastPara.Range.Find.Execute(
    "Prompt: Create a feature image for Aztec Warrior that showcases the game's exciting theme and features. The image should be in cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be standing in a jungle setting, with the Aztec temple in the background. The image should also include elements of the game's features, such as the Aztec Sun symbol and the expanding symbols during the free spin feature. The overall design should be bright, colorful, and eye-catching to attract players to the game.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Experience free spins and expanding symbols in Aztec Warrior, a 5-reel, 3-row, 10-payline slot game with an RTP of 96.50% - highly recommended for experienced gamblers.",
    2)
